$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "26.657.21"
$ws.Cells.Item(2, 5).Value = "  -0.49%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "1.631.27"
$ws.Cells.Item(3, 5).Value = "  -1.04%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  -0.19%  "

# Row 5
$ws.Cells.Item(5, 4).Value = "'217.90"
$ws.Cells.Item(5, 5).Value = "  +0.63%  "

# Row 6
$ws.Cells.Item(6, 4).Value = "'0.497"
$ws.Cells.Item(6, 5).Value = "  -1.74%  "

# Row 7
$ws.Cells.Item(7, 5).Value = "  -0.12%  "

# Row 8
$ws.Cells.Item(8, 4).Value = "'0.248"
$ws.Cells.Item(8, 5).Value = "  -1.44%  "

# Row 9
$ws.Cells.Item(9, 5).Value = "  -1.10%  "

# Row 10
$ws.Cells.Item(10, 4).Value = "'18.95"
$ws.Cells.Item(10, 5).Value = "  -1.51%  "

# Row 11
$ws.Cells.Item(11, 4).Value = "'0.0842"
$ws.Cells.Item(11, 5).Value = "  -0.18%  "

# Row 12
$ws.Cells.Item(12, 4).Value = "1.859.92"

# Row 13
$ws.Cells.Item(13, 4).Value = "1.629.41"
$ws.Cells.Item(13, 5).Value = "  -1.27%  "

# Row 14
$ws.Cells.Item(14, 4).Value = "'4.10"
$ws.Cells.Item(14, 5).Value = "  -2.36%  "

# Row 15
$ws.Cells.Item(15, 5).Value = "  -2.28%  "

# Row 16
$ws.Cells.Item(16, 4).Value = "'63.90"
$ws.Cells.Item(16, 5).Value = "  -2.37%  "

# Row 17
$ws.Cells.Item(17, 4).Value = "26.663.62"
$ws.Cells.Item(17, 5).Value = "  -0.50%  "

# Row 18
$ws.Cells.Item(18, 4).Value = "0.0₃0720"
$ws.Cells.Item(18, 5).Value = "  -3.24%  "

# Row 19
$ws.Cells.Item(19, 4).Value = "'210.98"
$ws.Cells.Item(19, 5).Value = "  -3.17%  "

# Row 20
$ws.Cells.Item(20, 5).Value = "  -0.10%  "

# Row 21
$ws.Cells.Item(21, 5).Value = "  -1.75%  "

# Row 22
$ws.Cells.Item(22, 5).Value = "  -7.06%  "

# Row 23
$ws.Cells.Item(23, 5).Value = "  -2.48%  "

# Row 24
$ws.Cells.Item(24, 5).Value = "  -3.27%  "

# Row 25
$ws.Cells.Item(25, 4).Value = "'146.55"
$ws.Cells.Item(25, 5).Value = "  +0.61%  "

# Row 26
$ws.Cells.Item(26, 5).Value = "  -0.07%  "

# Row 27
$ws.Cells.Item(27, 5).Value = "  -2.75%  "

# Row 28
$ws.Cells.Item(28, 5).Value = "  -2.78%  "

# Row 29
$ws.Cells.Item(29, 4).Value = "'15.48"
$ws.Cells.Item(29, 5).Value = "  -2.16%  "

# Row 30
$ws.Cells.Item(30, 5).Value = "  -3.60%  "

# Row 31
$ws.Cells.Item(31, 5).Value = "  +0.80%  "

# Row 32
$ws.Cells.Item(32, 5).Value = "  +0.35%  "

# Row 33
$ws.Cells.Item(33, 5).Value = "  -2.65%  "

# Row 34
$ws.Cells.Item(34, 4).Value = "1.257.89"
$ws.Cells.Item(34, 5).Value = "  -1.96%  "

# Row 35
$ws.Cells.Item(35, 2).Value = "LidoDAOToken"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Cells.Item(35, 4).Value = "'1.51"
$ws.Cells.Item(35, 5).Value = "  -2.46%  "

# Row 36
$ws.Cells.Item(36, 2).Value = "HuobiToken"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Cells.Item(36, 4).Value = "'2.44"
$ws.Cells.Item(36, 5).Value = "  +0.08%  "

# Row 37
$ws.Cells.Item(37, 5).Value = "  -3.34%  "

# Row 38
$ws.Cells.Item(38, 5).Value = "  -3.33%  "

# Row 39
$ws.Cells.Item(39, 5).Value = "  -0.15%  "

# Row 40
$ws.Cells.Item(40, 4).Value = "'0.800"
$ws.Cells.Item(40, 5).Value = "  -3.92%  "

# Row 41
$ws.Cells.Item(41, 5).Value = "  -2.58%  "

# Row 42
$ws.Cells.Item(42, 4).Value = "1.770.57"
$ws.Cells.Item(42, 5).Value = "  -1.53%  "

# Row 43
$ws.Cells.Item(43, 5).Value = "  -4.60%  "

# Row 44
$ws.Cells.Item(44, 5).Value = "  -3.99%  "

# Row 45
$ws.Cells.Item(45, 4).Value = "'90.91"

# Row 46
$ws.Cells.Item(46, 4).Value = "'59.68"
$ws.Cells.Item(46, 5).Value = "  -0.13%  "

# Row 47
$ws.Cells.Item(47, 5).Value = "  -2.78%  "

# Row 48
$ws.Cells.Item(48, 5).Value = "  -0.15%  "

# Row 49
$ws.Cells.Item(49, 5).Value = "  -0.11%  "

# Row 50
$ws.Cells.Item(50, 2).Value = "EnergySwap"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(50, 4).Value = "'7.47"
$ws.Cells.Item(50, 5).Value = "  -3.88%  "

# Row 51
$ws.Cells.Item(51, 5).Value = "  -2.83%  "
